# Actualizado al 18 de agosto 2020
# Appends the 17-Aug-2020 (44060) and 18-Aug-2020 (44061) rows to every
# per-department sheet, matching the upstream COVID-19 Bolivia workbook
# commit "Actulizado al 18 de agosto 2020".

$wb = $excel.ActiveWorkbook

# sheet (tab name) -> department label, new-row values for 44060 / 44061
# values are [Confirmados, Decesos, Recuperados]
$data = @(
    @{ Sheet = "bn"; Depto = "Beni";        D1 = @(14, 0, 3);   D2 = @(51, 1, 4) },
    @{ Sheet = "cb"; Depto = "Cochabamba";  D1 = @(52, 7, 323); D2 = @(80, 8, 320) },
    @{ Sheet = "ch"; Depto = "Chuquisaca";  D1 = @(177, 8, 25); D2 = @(162, 9, 72) },
    @{ Sheet = "lp"; Depto = "La Paz";      D1 = @(276, 15, 69); D2 = @(243, 8, 149) },
    @{ Sheet = "or"; Depto = "Oruro";       D1 = @(14, 1, 20);  D2 = @(16, 1, 9) },
    @{ Sheet = "pn"; Depto = "Pando";       D1 = @(9, 0, 4);    D2 = @(53, 0, 0) },
    @{ Sheet = "pt"; Depto = "Potosí";      D1 = @(138, 14, 1); D2 = @(861, 4, 91) },
    @{ Sheet = "sc"; Depto = "Santa Cruz";  D1 = @(95, 16, 466); D2 = @(69, 15, 479) },
    @{ Sheet = "tj"; Depto = "Tarija";      D1 = @(104, 4, 69); D2 = @(261, 3, 28) }
)

foreach ($entry in $data) {
    $ws = $wb.Worksheets.Item($entry.Sheet)

    # Rows 162 and 163 don't exist yet (or only exist as empty placeholder
    # cells) in the source file. Copy the formatting of the last populated
    # row (161) down so the new rows pick up the same cell styles (date
    # number format on A, borders/fill on A:E, etc.) that Excel would keep
    # when a user fills the row below.
    $ws.Range("A161:E161").Copy() | Out-Null
    $ws.Range("A162:E162").PasteSpecial(-4122) | Out-Null
    $ws.Range("A163:E163").PasteSpecial(-4122) | Out-Null

    $ws.Range("A162").Value = 44060
    $ws.Range("B162").Value = $entry.Depto
    $ws.Range("C162").Value = $entry.D1[0]
    $ws.Range("D162").Value = $entry.D1[1]
    $ws.Range("E162").Value = $entry.D1[2]

    $ws.Range("A163").Value = 44061
    $ws.Range("B163").Value = $entry.Depto
    $ws.Range("C163").Value = $entry.D2[0]
    $ws.Range("D163").Value = $entry.D2[1]
    $ws.Range("E163").Value = $entry.D2[2]
}

# "lp" (La Paz) gained two more blank trailer rows (164/165) carrying only
# the A/B column formatting, same shape as the blank placeholder row the
# other sheets had at row 162 before this edit.
$wsLp = $wb.Worksheets.Item("lp")
$wsLp.Range("A161:B161").Copy() | Out-Null
$wsLp.Range("A164:B164").PasteSpecial(-4122) | Out-Null
$wsLp.Range("A165:B165").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# Update each tab's remembered selection, matching what Excel would leave
# behind after entering the new rows on that sheet.
$wb.Worksheets.Item("bn").Range("E162").Select() | Out-Null
$wb.Worksheets.Item("cb").Range("C163").Select() | Out-Null
$wb.Worksheets.Item("ch").Range("E162").Select() | Out-Null
$wb.Worksheets.Item("lp").Range("E163").Select() | Out-Null
$wb.Worksheets.Item("or").Range("F162").Select() | Out-Null
$wb.Worksheets.Item("pn").Range("D163").Select() | Out-Null
$wb.Worksheets.Item("pt").Range("E163").Select() | Out-Null
$wb.Worksheets.Item("sc").Range("C163").Select() | Out-Null
$wb.Worksheets.Item("tj").Range("E164").Select() | Out-Null

# The active tab moved from "tj" (last sheet) to "bn" (first sheet).
$wb.Worksheets.Item("bn").Activate()
$wb.Worksheets.Item("bn").Range("E162").Select() | Out-Null
